# Apply the update that appends four new daily rows (166-169) to the
# Arequipa COVID data sheet and normalizes three stray cell styles
# (E161, J161, C162) that were using the wrong font style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 0: a few new cells (B168, P169, Q169) must end up using the same
# "odd" cell style (style index 13 in the original file) that currently
# sits on E161 / J161 / C162. Capture that style now, before we correct
# those three cells below, by stashing a copy of it on a scratch cell far
# away from any real data. We clean the scratch cell up at the end.
# ---------------------------------------------------------------------
$scratch = "A1000"
$ws.Range("E161").Copy()
$ws.Range($scratch).PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Step 1: fix the three mis-styled existing cells so they match the
# rest of their row (copy format from an adjacent, correctly styled
# cell in the same row).
# ---------------------------------------------------------------------
$ws.Range("D161").Copy()
$ws.Range("E161").PasteSpecial($xlPasteFormats)

$ws.Range("I161").Copy()
$ws.Range("J161").PasteSpecial($xlPasteFormats)

$ws.Range("B162").Copy()
$ws.Range("C162").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Step 2: build out rows 166-169 with the correct formatting.
#   - Column A uses the date style already used by A163:A165.
#   - Columns B:Q use the plain style already used by B161:Q161
#     (which is now fully consistent after Step 1).
# ---------------------------------------------------------------------
$ws.Range("A163").Copy()
$ws.Range("A166:A169").PasteSpecial($xlPasteFormats)

$ws.Range("B161:Q161").Copy()
$ws.Range("B166:Q166").PasteSpecial($xlPasteFormats)
$ws.Range("B167:Q167").PasteSpecial($xlPasteFormats)
$ws.Range("B168:Q168").PasteSpecial($xlPasteFormats)
$ws.Range("B169:Q169").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Step 3: re-apply the stashed "odd" style to the three cells that need
# it (B168, P169, Q169), then discard the scratch cell entirely.
# ---------------------------------------------------------------------
$ws.Range($scratch).Copy()
$ws.Range("B168").PasteSpecial($xlPasteFormats)
$ws.Range("P169").PasteSpecial($xlPasteFormats)
$ws.Range("Q169").PasteSpecial($xlPasteFormats)

$ws.Range($scratch).Clear()

# ---------------------------------------------------------------------
# Step 4: fill in the values for the new rows.
# ---------------------------------------------------------------------
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

$row166 = @(44116,828457,132154,694721,1582,2158,121391,266,54,60,0,9,817,1170,38,123,10)
$row167 = @(44117,831368,132417,697083,1868,2159,122007,246,61,56,0,9,818,1170,38,123,10)
$row168 = @(44118,834358,132754,700242,1362,2161,122974,240,56,52,0,7,818,1172,38,123,10)
$row169 = @(44119,836958,132964,702645,1349,2168,123918,240,56,52,0,7,822,1175,38,123,10)

$rows = @{ 166 = $row166; 167 = $row167; 168 = $row168; 169 = $row169 }

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

Write-Output "Applied Arequipa COVID data update (rows 166-169 + style fixes)."
